# Applies the "value-calculator-logic" number tweaks described in the
# commit diff. Two slides (11 and 12) each contain a small set of
# textboxes (several nested inside group shapes) whose numbers must be
# swapped out while leaving every other bit of markup (runs, paragraphs,
# the existing <a:br/> soft line-breaks, run formatting, shape sizing,
# etc.) untouched.
#
# Two subtleties of the object model are handled explicitly:
#
#  1. Overwriting TextFrame.TextRange.Text wholesale collapses the
#     existing multi-run / <a:br/> structure into fresh paragraphs and
#     stamps a brand new <a:rPr lang="en-US"/> onto the run. Instead,
#     every edit below targets the exact Characters(start, length) span
#     that corresponds to one whole original <a:r> run and rewrites just
#     that span with the full new run text - that keeps the run a single
#     run, leaves its (missing) rPr alone, and leaves sibling <a:br/>
#     elements untouched.
#
#  2. These textboxes all use <a:spAutoFit/>, so committing a text edit
#     makes PowerPoint re-lay-out and resize the shape (Shape.Height
#     changes). The source deck's own text-to-box-size fit never moved,
#     so each edit below snapshots Shape.Height beforehand and reapplies
#     it afterward to keep the shape's size as it was.

$p = $ppt.ActivePresentation

function Set-RunSpan($textRange, $start, $len, $newText) {
    $textRange.Characters($start, $len).Text = $newText
}

# Replace the entire text of a shape that holds exactly one run/one
# paragraph, addressing the run's full span so no <a:rPr> gets stamped in,
# and restoring the autofit-recalculated height afterward.
function Set-WholeShapeText($shape, $newText) {
    $origHeight = $shape.Height
    $tr = $shape.TextFrame.TextRange
    Set-RunSpan $tr 1 $tr.Text.Length $newText
    $shape.Height = $origHeight
}

# ---------------------------------------------------------------------
# Slide 11
# ---------------------------------------------------------------------
$s11 = $p.Slides.Item(11)

# "■契約更新対象者数 : 111名" / "■店長の時給 : 2,222円" / "■契約更新の回数 : 年3回"
# (one paragraph, three runs separated by <a:br/>) -> update runs 1 & 2.
$shp = $s11.Shapes.Item(5).GroupItems.Item(2)
$h = $shp.Height
$tr = $shp.TextFrame.TextRange
Set-RunSpan $tr 18 15 "■店長の時給 : 12円"
Set-RunSpan $tr 1 16 "■契約更新対象者数 : 11名"
$shp.Height = $h

# "111人" -> "11人"
Set-WholeShapeText $s11.Shapes.Item(6).GroupItems.Item(2) "11人"

# "124分" -> "65分"
Set-WholeShapeText $s11.Shapes.Item(7).GroupItems.Item(2) "65分"

# "41,292分" -> "2,145分"
Set-WholeShapeText $s11.Shapes.Item(8).GroupItems.Item(2) "2,145分"

# "年間で" + <a:br/> + "688時間" -> "年間で" + <a:br/> + "36時間"
$shp = $s11.Shapes.Item(11).GroupItems.Item(2)
$h = $shp.Height
$tr = $shp.TextFrame.TextRange
Set-RunSpan $tr 5 5 "36時間"
$shp.Height = $h

# "111人" -> "11人"
Set-WholeShapeText $s11.Shapes.Item(15).GroupItems.Item(2) "11人"

# "999分" -> "99分"
Set-WholeShapeText $s11.Shapes.Item(17).GroupItems.Item(2) "99分"

# "年間で" + <a:br/> + "17時間" -> "年間で" + <a:br/> + "2時間"
$shp = $s11.Shapes.Item(20).GroupItems.Item(2)
$h = $shp.Height
$tr = $shp.TextFrame.TextRange
Set-RunSpan $tr 5 5 "2時間"
$shp.Height = $h

# ---------------------------------------------------------------------
# Slide 12
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)

# "■契約更新対象者数 : 111名" / "■店長の時給 : 2,222円" / "■契約更新の回数 : 年3回"
$shp = $s12.Shapes.Item(5).GroupItems.Item(2)
$h = $shp.Height
$tr = $shp.TextFrame.TextRange
Set-RunSpan $tr 18 15 "■店長の時給 : 12円"
Set-RunSpan $tr 1 16 "■契約更新対象者数 : 11名"
$shp.Height = $h

# "688時間" -> "36時間"
Set-WholeShapeText $s12.Shapes.Item(6).GroupItems.Item(2) "36時間"

# "2222円" -> "12円"
Set-WholeShapeText $s12.Shapes.Item(7).GroupItems.Item(2) "12円"

# "4,586,208円" -> "1,296円"
Set-WholeShapeText $s12.Shapes.Item(8).GroupItems.Item(2) "1,296円"

# "17時間" -> "2時間"
Set-WholeShapeText $s12.Shapes.Item(15).GroupItems.Item(2) "2時間"

# "102,000円" -> "12,000円"
Set-WholeShapeText $s12.Shapes.Item(17).GroupItems.Item(2) "12,000円"

# "151,950円" -> "16,950円"
Set-WholeShapeText $s12.Shapes.Item(22).GroupItems.Item(2) "16,950円"

# "111人" -> "11人"
Set-WholeShapeText $s12.Shapes.Item(24).GroupItems.Item(2) "11人"

# "49,950円" -> "4,950円"
Set-WholeShapeText $s12.Shapes.Item(26).GroupItems.Item(2) "4,950円"
